$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overall status -> Complete for rows 10, 11, 12, 14
$ws.Range("H10").Value = "Complete"
$ws.Range("H11").Value = "Complete"
$ws.Range("H12").Value = "Complete"
$ws.Range("H14").Value = "Complete"

# Clear notes column (J) for rows 6, 7, 13, 15
$ws.Range("J6").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("J15").ClearContents()

# Assigned To (I) for rows 16-19
$ws.Range("I16").Value = "Austin"
$ws.Range("I17").Value = "Austin"
$ws.Range("I18").Value = "Ausin"
$ws.Range("I19").Value = "Richard"

# Highlight row 19 (A:K) with the same fill used elsewhere in the sheet
$ws.Range("A19:K19").Interior.ThemeColor = 2
$ws.Range("A19:K19").Interior.TintAndShade = -0.249977111117893

# Update selection
$ws.Range("E32").Select()

Write-Output "done"
